# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" (fund holdings detail) positioned
#    between "总计" and "2021-Q4".
# 2. Insert a new summary row into "总计" for 2022-Q3, pushing the
#    existing 2021-Q4 summary row down.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q4    = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. New "2022-Q3" sheet with per-fund holding detail, inserted right
#    before the existing "2021-Q4" sheet.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($q4)
$q3.Name = "2022-Q3"

# Header row
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"
$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.Item(7).LineStyle = 1
$headerRange.Borders.Item(8).LineStyle = 1
$headerRange.Borders.Item(9).LineStyle = 1
$headerRange.Borders.Item(10).LineStyle = 1
$headerRange.Borders.Item(11).LineStyle = 1
$headerRange.Borders.Item(12).LineStyle = 1

# Fund rows (code/name/size/position/weight/value are stored as text,
# matching the source data which keeps the original string formatting;
# the index column A and rank column H are numeric).
$rows = @(
    @(0,"000006","西部利得量化成长混合A","13.82","86.11","1.05","0.1451",8),
    @(1,"011228","西部利得量化成长混合C","1.67","86.11","1.05","0.0175",8),
    @(2,"010703","财通智选消费股票A","0.48","93.51","3.59","0.0172",4),
    @(3,"010704","财通智选消费股票C","0.46","93.51","3.59","0.0165",4),
    @(4,"014214","光大保德信核心资产混合A","0.32","88.01","3.57","0.0114",10),
    @(5,"014462","光大保德信汇佳混合A","0.33","43.38","3.23","0.0107",1),
    @(6,"014215","光大保德信核心资产混合C","0.03","88.01","3.57","0.0011",10),
    @(7,"014463","光大保德信汇佳混合C","0.03","43.38","3.23","0.0010",1)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q3.Cells.Item($r,1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.Item(7).LineStyle = 1
    $idxCell.Borders.Item(8).LineStyle = 1
    $idxCell.Borders.Item(9).LineStyle = 1
    $idxCell.Borders.Item(10).LineStyle = 1

    $q3.Cells.Item($r,2).Value = "'" + $row[1]
    $q3.Cells.Item($r,3).Value = $row[2]
    $q3.Cells.Item($r,4).Value = "'" + $row[3]
    $q3.Cells.Item($r,5).Value = "'" + $row[4]
    $q3.Cells.Item($r,6).Value = "'" + $row[5]
    $q3.Cells.Item($r,7).Value = "'" + $row[6]
    $q3.Cells.Item($r,8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. "总计" sheet: insert a new row 2 for 2022-Q3, shifting the existing
#    2021-Q4 summary row down to row 3.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").Style = "Normal"

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 8
$total.Cells.Item(2,4).Value = 0.22

$a2 = $total.Cells.Item(2,1)
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.Item(7).LineStyle = 1
$a2.Borders.Item(8).LineStyle = 1
$a2.Borders.Item(9).LineStyle = 1
$a2.Borders.Item(10).LineStyle = 1

# The shifted-down 2021-Q4 row keeps its own running index, now 1.
$total.Cells.Item(3,1).Value = 1
